$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.415.46"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "1.820.69"
$ws.Range("E3").Value = "  +1.82%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "

$ws.Range("E6").Value = "  +1.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.11"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.66%  "

$ws.Range("E9").Value = "  +0.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0694"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0949"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("D12").Value = "2.078.71"
$ws.Range("E12").Value = "  +1.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.14%  "

$ws.Range("D14").Value = "1.810.42"
$ws.Range("E14").Value = "  -0.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.649"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.70%  "

$ws.Range("D16").Value = "34.445.91"
$ws.Range("E16").Value = "  +0.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.28%  "

$ws.Range("D19").Value = "0.0₃0801"
$ws.Range("E19").Value = "  -0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.40%  "

$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "167.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.77%  "

$ws.Range("E25").Value = "  +1.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.49%  "

$ws.Range("E28").Value = "  +1.25%  "

$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0530"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.60%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.16%  "

$ws.Range("E35").Value = "  +1.10%  "

$ws.Range("D36").Value = "1.414.02"
$ws.Range("E36").Value = "  -2.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.679"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.42%  "

$ws.Range("B39").Value = "Aave"
$ws.Range("C39").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "86.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.82%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0191"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.961"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.37%  "

$ws.Range("E43").Value = "  +1.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0522"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.53%  "

$ws.Range("E46").Value = "  +3.02%  "

$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("D48").Value = "1.981.15"
$ws.Range("E48").Value = "  +1.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").Value = "0.0₆0129"
$ws.Range("E51").Value = "  +1.61%  "
